$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.950.76"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.818.16"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'310.32"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "'0.07346"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "'0.8726"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'20.27"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "1.824.66"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "'5.405"
$ws.Range("D14").Value = "'0.07116"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'6.508"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'91.40"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'0.000008709"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'14.65"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "26.968.42"
$ws.Range("D22").Value = "'5.279"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "2.046.24"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'1.894"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "'150.98"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'2.151"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'5.242"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "'116.94"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").Value = "'0.08889"
$ws.Range("D32").Value = "'0.7582"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "'4.507"
$ws.Range("D35").Value = "'2.915"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D37").Value = "'1.096"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'0.05308"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'0.01945"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'2.967"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "'2.380"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("D42").Value = "'0.5296"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'7.173"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "'0.1654"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "'8.439"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "'103.44"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'1.660"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "'0.06302"
$ws.Range("E51").Value = "  +0.34%  "

Write-Host "Updated cryptos list values"
